# Update countries & provincias Spain
# Updates the "last updated" timestamp and refreshes COVID-19 stats for a
# handful of countries on the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" footer text (row 1, column A) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 08:17"

# --- New data values per country: Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes ---
$updates = @(
    @{ Country = "Estados Unidos"; B = 2388225; C = 72;   D = 1003062; E = 1262552; F = 0; G = 1;   H = 122611 },
    @{ Country = "India";          B = 440685;  C = 235;  D = 248190;  E = 178480;  F = 0; G = 0;   H = 14015 },
    @{ Country = "Pakistan";       B = 185034;  C = 3946; D = 73471;   E = 107868;  F = 0; G = 105; H = 3695 },
    @{ Country = "Afganistan";     B = 29481;   C = 324;  D = 9260;    E = 19603;   F = 0; G = 20;  H = 618 },
    @{ Country = "Uzbekistan";     B = 6500;    C = 39;   D = 4450;    E = 2031;    F = 0; G = 0;   H = 19 },
    @{ Country = "Haiti";          B = 5211;    C = 0;    D = 338;     E = 4785;    F = 0; G = 0;   H = 88 },
    @{ Country = "El Salvador";    B = 4808;    C = 0;    D = 2695;    E = 2006;    F = 0; G = 0;   H = 107 },
    @{ Country = "Hungria";        B = 4107;    C = 5;    D = 2600;    E = 934;     F = 0; G = 1;   H = 573 },
    @{ Country = "Kirguistan";     B = 3519;    C = 163;  D = 2054;    E = 1424;    F = 0; G = 1;   H = 41 },
    @{ Country = "Tailandia";      B = 3156;    C = 5;    D = 3023;    E = 75;      F = 0; G = 0;   H = 58 },
    @{ Country = "Camboya";        B = 130;     C = 1;    D = 127;     E = 3;       F = 0; G = 0;   H = 0 },
    @{ Country = "Butan";          B = 69;      C = 1;    D = 32;      E = 37;      F = 0; G = 0;   H = 0 }
)

foreach ($u in $updates) {
    $found = $ws.Columns.Item(1).Find($u.Country)
    $row = $found.Row
    $ws.Cells.Item($row, 2).Value = $u.B
    $ws.Cells.Item($row, 3).Value = $u.C
    $ws.Cells.Item($row, 4).Value = $u.D
    $ws.Cells.Item($row, 5).Value = $u.E
    $ws.Cells.Item($row, 6).Value = $u.F
    $ws.Cells.Item($row, 7).Value = $u.G
    $ws.Cells.Item($row, 8).Value = $u.H
}
